$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value2 = 0.2088353413654618
$ws.Cells.Item(2, 3).Value2 = 0.5341365461847389
$ws.Cells.Item(2, 10).Value2 = 0.01606425702811245
$ws.Cells.Item(2, 16).Value2 = 0.1526104417670683
$ws.Cells.Item(2, 19).Value2 = 0.08835341365461848
$ws.Cells.Item(3, 3).Value2 = 0.01481481481481482
$ws.Cells.Item(3, 10).Value2 = 0.05925925925925926
$ws.Cells.Item(3, 16).Value2 = 0.7185185185185186
$ws.Cells.Item(3, 19).Value2 = 0.2074074074074074
$ws.Cells.Item(4, 10).Value2 = 0.05714285714285714
$ws.Cells.Item(4, 16).Value2 = 0.6285714285714286
$ws.Cells.Item(4, 19).Value2 = 0.3142857142857143
$ws.Cells.Item(6, 2).Value2 = 0.05555555555555555
$ws.Cells.Item(6, 4).Value2 = 0.01851851851851852
$ws.Cells.Item(6, 6).Value2 = 0.01851851851851852
$ws.Cells.Item(6, 10).Value2 = 0.2777777777777778
$ws.Cells.Item(6, 15).Value2 = 0.01851851851851852
$ws.Cells.Item(6, 17).Value2 = 0.1342592592592593
$ws.Cells.Item(6, 18).Value2 = 0.03703703703703703
$ws.Cells.Item(6, 19).Value2 = 0.4398148148148148
$ws.Cells.Item(7, 2).Value2 = 0.05752212389380531
$ws.Cells.Item(7, 4).Value2 = 0.02212389380530973
$ws.Cells.Item(7, 5).Value2 = 0.004424778761061947
$ws.Cells.Item(7, 6).Value2 = 0.07964601769911504
$ws.Cells.Item(7, 10).Value2 = 0.1415929203539823
$ws.Cells.Item(7, 15).Value2 = 0.01769911504424779
$ws.Cells.Item(7, 17).Value2 = 0.168141592920354
$ws.Cells.Item(7, 18).Value2 = 0.06637168141592921
$ws.Cells.Item(7, 19).Value2 = 0.4424778761061947
$ws.Cells.Item(8, 2).Value2 = 0.06813186813186813
$ws.Cells.Item(8, 4).Value2 = 0.01978021978021978
$ws.Cells.Item(8, 5).Value2 = 0.002197802197802198
$ws.Cells.Item(8, 6).Value2 = 0.06373626373626373
$ws.Cells.Item(8, 10).Value2 = 0.0945054945054945
$ws.Cells.Item(8, 15).Value2 = 0.01098901098901099
$ws.Cells.Item(8, 17).Value2 = 0.2065934065934066
$ws.Cells.Item(8, 18).Value2 = 0.08131868131868132
$ws.Cells.Item(8, 19).Value2 = 0.4527472527472527
$ws.Cells.Item(9, 2).Value2 = 0.07253886010362694
$ws.Cells.Item(9, 4).Value2 = 0.0155440414507772
$ws.Cells.Item(9, 6).Value2 = 0.05699481865284974
$ws.Cells.Item(9, 10).Value2 = 0.1036269430051813
$ws.Cells.Item(9, 15).Value2 = 0.02590673575129534
$ws.Cells.Item(9, 17).Value2 = 0.1450777202072539
$ws.Cells.Item(9, 18).Value2 = 0.07253886010362694
$ws.Cells.Item(9, 19).Value2 = 0.5077720207253886
$ws.Cells.Item(10, 2).Value2 = 0.09948979591836735
$ws.Cells.Item(10, 4).Value2 = 0.01530612244897959
$ws.Cells.Item(10, 5).Value2 = 0.001700680272108843
$ws.Cells.Item(10, 6).Value2 = 0.06802721088435375
$ws.Cells.Item(10, 10).Value2 = 0.1020408163265306
$ws.Cells.Item(10, 15).Value2 = 0.02040816326530612
$ws.Cells.Item(10, 17).Value2 = 0.2346938775510204
$ws.Cells.Item(10, 18).Value2 = 0.06462585034013606
$ws.Cells.Item(10, 19).Value2 = 0.3937074829931973
$ws.Cells.Item(11, 7).Value2 = 0.1554959785522788
$ws.Cells.Item(11, 10).Value2 = 0.08579088471849866
$ws.Cells.Item(11, 11).Value2 = 0.2359249329758713
$ws.Cells.Item(11, 12).Value2 = 0.5013404825737265
$ws.Cells.Item(11, 19).Value2 = 0.02144772117962467
$ws.Cells.Item(12, 7).Value2 = 0.7157894736842105
$ws.Cells.Item(12, 10).Value2 = 0.2526315789473684
$ws.Cells.Item(12, 11).Value2 = 0.01052631578947368
$ws.Cells.Item(12, 12).Value2 = 0.005263157894736842
$ws.Cells.Item(12, 19).Value2 = 0.01578947368421053
$ws.Cells.Item(13, 6).Value2 = 0.01754385964912281
$ws.Cells.Item(13, 10).Value2 = 0.2456140350877193
$ws.Cells.Item(13, 19).Value2 = 0.07017543859649122
$ws.Cells.Item(15, 6).Value2 = 0.0101010101010101
$ws.Cells.Item(15, 8).Value2 = 0.1868686868686869
$ws.Cells.Item(15, 9).Value2 = 0.0505050505050505
$ws.Cells.Item(15, 10).Value2 = 0.303030303030303
$ws.Cells.Item(15, 11).Value2 = 0.06060606060606061
$ws.Cells.Item(15, 13).Value2 = 0.0202020202020202
$ws.Cells.Item(15, 14).Value2 = 0.005050505050505051
$ws.Cells.Item(15, 15).Value2 = 0.08080808080808081
$ws.Cells.Item(15, 19).Value2 = 0.2828282828282828
$ws.Cells.Item(16, 6).Value2 = 0.02
$ws.Cells.Item(16, 8).Value2 = 0.1333333333333333
$ws.Cells.Item(16, 9).Value2 = 0.07333333333333333
$ws.Cells.Item(16, 10).Value2 = 0.4266666666666667
$ws.Cells.Item(16, 11).Value2 = 0.1
$ws.Cells.Item(16, 15).Value2 = 0.07333333333333333
$ws.Cells.Item(16, 19).Value2 = 0.1733333333333333
$ws.Cells.Item(17, 6).Value2 = 0.02844638949671772
$ws.Cells.Item(17, 8).Value2 = 0.1553610503282276
$ws.Cells.Item(17, 9).Value2 = 0.07439824945295405
$ws.Cells.Item(17, 10).Value2 = 0.387308533916849
$ws.Cells.Item(17, 11).Value2 = 0.1312910284463895
$ws.Cells.Item(17, 13).Value2 = 0.03063457330415755
$ws.Cells.Item(17, 14).Value2 = 0.00437636761487965
$ws.Cells.Item(17, 15).Value2 = 0.04595185995623632
$ws.Cells.Item(17, 19).Value2 = 0.1422319474835886
$ws.Cells.Item(18, 8).Value2 = 0.2416107382550336
$ws.Cells.Item(18, 9).Value2 = 0.1006711409395973
$ws.Cells.Item(18, 10).Value2 = 0.3691275167785235
$ws.Cells.Item(18, 11).Value2 = 0.0738255033557047
$ws.Cells.Item(18, 13).Value2 = 0.02684563758389262
$ws.Cells.Item(18, 15).Value2 = 0.02013422818791946
$ws.Cells.Item(18, 19).Value2 = 0.1677852348993289
$ws.Cells.Item(19, 6).Value2 = 0.02352941176470588
$ws.Cells.Item(19, 8).Value2 = 0.2161764705882353
$ws.Cells.Item(19, 9).Value2 = 0.09117647058823529
$ws.Cells.Item(19, 10).Value2 = 0.3345588235294117
$ws.Cells.Item(19, 11).Value2 = 0.1330882352941176
$ws.Cells.Item(19, 13).Value2 = 0.025
$ws.Cells.Item(19, 14).Value2 = 0.002205882352941176
$ws.Cells.Item(19, 15).Value2 = 0.05882352941176471
$ws.Cells.Item(19, 19).Value2 = 0.1154411764705882
